# Applies the commit's change:
#  - Insert a new "Player Info" worksheet as the first sheet in the workbook,
#    populated with the player's ID/NAME/BATTING_HAND/BOWL_STYLE.
#  - On "ODI Batting" and "ODI Bowling", rename the MATCH_CARD_LINK column to
#    MATCH_CODE and replace the full scorecard URL value with just the
#    numeric match code extracted from it.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet before "ODI Batting" ---
# NOTE: worksheet variables captured before a sheet insertion can become
# stale (they track a position, not a stable identity), so re-fetch sheets
# by name after any call that changes the sheet collection.
$battingWsForInsert = $wb.Worksheets.Item("ODI Batting")
$infoWs = $wb.Worksheets.Add($battingWsForInsert)
$infoWs.Name = "Player Info"

# --- 2. Populate "Player Info" header row + data row ---
$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $infoHeaders.Length; $i++) {
    $cell = $infoWs.Cells.Item(1, $i + 1)
    $cell.Value = $infoHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Leading apostrophe forces the numeric-looking ID to be stored as text,
# matching the other sheets where similar values are text.
$infoWs.Cells.Item(2, 1).Value = "'6624"
$infoWs.Cells.Item(2, 2).Value = "Oliver Davidson"
$infoWs.Cells.Item(2, 3).Value = "Left Handed"
$infoWs.Cells.Item(2, 4).Value = "Left Arm Orthodox"

# --- 3. Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D) ---
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Cells.Item(1, 4).Value = "MATCH_CODE"
$battingWs.Cells.Item(2, 4).Value = "'4581"

# --- 4. Update "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B) ---
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$bowlingWs.Cells.Item(1, 2).Value = "MATCH_CODE"
$bowlingWs.Cells.Item(2, 2).Value = "'4581"
